# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet as the last sheet ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the exact same header / date formatting already used on the
# "Weekly Quantity" sheet so no duplicate style entries are created.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122) # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122) # xlPasteFormats

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$wsForecast.Cells.Item(2,1).Value = 45578.99999999999
$wsForecast.Cells.Item(2,2).Value = 186
$wsForecast.Cells.Item(2,3).Value = 110.6953038070558
$wsForecast.Cells.Item(2,4).Value = 266.8887152075238

$wsForecast.Cells.Item(3,1).Value = 45592.99999999999
$wsForecast.Cells.Item(3,2).Value = 104
$wsForecast.Cells.Item(3,3).Value = 35.54131260877912
$wsForecast.Cells.Item(3,4).Value = 180.9395897150247

$wsForecast.Cells.Item(4,1).Value = 45613.99999999999
$wsForecast.Cells.Item(4,2).Value = 0
$wsForecast.Cells.Item(4,3).Value = -95.13524414842824
$wsForecast.Cells.Item(4,4).Value = 56.5066942586831

$wsForecast.Cells.Item(5,1).Value = 45620.99999999999
$wsForecast.Cells.Item(5,2).Value = 0
$wsForecast.Cells.Item(5,3).Value = -135.4230749981228
$wsForecast.Cells.Item(5,4).Value = 16.69113352466209

$wsForecast.Cells.Item(6,1).Value = 45627.99999999999
$wsForecast.Cells.Item(6,2).Value = 0
$wsForecast.Cells.Item(6,3).Value = -178.6862910865262
$wsForecast.Cells.Item(6,4).Value = -25.45823239031739

$wsForecast.Cells.Item(7,1).Value = 45634.99999999999
$wsForecast.Cells.Item(7,2).Value = 0
$wsForecast.Cells.Item(7,3).Value = -221.3177421842085
$wsForecast.Cells.Item(7,4).Value = -63.87392219662905

$wsForecast.Cells.Item(8,1).Value = 45641.99999999999
$wsForecast.Cells.Item(8,2).Value = 0
$wsForecast.Cells.Item(8,3).Value = -261.4357728144859
$wsForecast.Cells.Item(8,4).Value = -106.9783768464102

$wsForecast.Cells.Item(9,1).Value = 45648.99999999999
$wsForecast.Cells.Item(9,2).Value = 0
$wsForecast.Cells.Item(9,3).Value = -303.9545022179584
$wsForecast.Cells.Item(9,4).Value = -151.1082456976727

$wsForecast.Cells.Item(10,1).Value = 45655.99999999999
$wsForecast.Cells.Item(10,2).Value = 0
$wsForecast.Cells.Item(10,3).Value = -338.0910265466177
$wsForecast.Cells.Item(10,4).Value = -188.8939996347721

$wsForecast.Cells.Item(11,1).Value = 45662.99999999999
$wsForecast.Cells.Item(11,2).Value = 0
$wsForecast.Cells.Item(11,3).Value = -381.0533636786504
$wsForecast.Cells.Item(11,4).Value = -225.58879712624

$wsForecast.Cells.Item(12,1).Value = 45669.99999999999
$wsForecast.Cells.Item(12,2).Value = 0
$wsForecast.Cells.Item(12,3).Value = -419.3204304136408
$wsForecast.Cells.Item(12,4).Value = -275.5428875924437

# Restore the originally active sheet/tab (adding a sheet makes it active)
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select()

Write-Output "PO Forecast sheet created with forecast data"
